$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.206.22"
$ws.Range("E2").Value = "  +1.66%  "

$ws.Range("D3").Value = "2.353.45"
$ws.Range("E3").Value = "  -1.22%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.678"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.44"
$ws.Range("E6").Value = "  +3.57%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.14"
$ws.Range("E7").Value = "  +6.49%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +23.24%  "

$ws.Range("E10").Value = "  +5.44%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "31.65"
$ws.Range("E11").Value = "  +20.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.39"
$ws.Range("E12").Value = "  +18.91%  "

$ws.Range("E13").Value = "  +2.46%  "

$ws.Range("D14").Value = "2.706.30"
$ws.Range("E14").Value = "  -1.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.81"
$ws.Range("E15").Value = "  +7.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.906"
$ws.Range("E16").Value = "  +6.68%  "

$ws.Range("D17").Value = "2.347.71"
$ws.Range("E17").Value = "  -1.09%  "

$ws.Range("D18").Value = "44.256.99"
$ws.Range("E18").Value = "  +1.72%  "

$ws.Range("E19").Value = "  +4.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.67"
$ws.Range("E20").Value = "  +5.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "78.00"
$ws.Range("E21").Value = "  +5.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "255.52"
$ws.Range("E22").Value = "  +1.83%  "

$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("E24").Value = "  -4.92%  "

$ws.Range("E25").Value = "  +4.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.70"
$ws.Range("E26").Value = "  +7.67%  "

$ws.Range("E27").Value = "  +3.68%  "

$ws.Range("E28").Value = "  -2.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "174.83"
$ws.Range("E29").Value = "  +1.25%  "

$ws.Range("E30").Value = "  +3.03%  "

$ws.Range("E31").Value = "  +3.60%  "

$ws.Range("E32").Value = "  +4.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.38"
$ws.Range("E33").Value = "  +8.40%  "

$ws.Range("E34").Value = "  +10.48%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.34"
$ws.Range("E35").Value = "  +5.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.87"
$ws.Range("E36").Value = "  +6.90%  "

$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.57"
$ws.Range("E37").Value = "  -0.29%  "

$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.44"
$ws.Range("E38").Value = "  -0.30%  "

$ws.Range("E39").Value = "  +7.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.32"
$ws.Range("E40").Value = "  +4.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.97"
$ws.Range("E41").Value = "  +0.71%  "

$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("E43").Value = "  +3.84%  "

$ws.Range("E44").Value = "  +5.47%  "

$ws.Range("E45").Value = "  +13.83%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.47"
$ws.Range("E46").Value = "  +10.56%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.61"
$ws.Range("E47").Value = "  +1.58%  "

$ws.Range("E48").Value = "  -1.96%  "

$ws.Range("E49").Value = "  -0.92%  "

$ws.Range("D50").Value = "1.450.94"
$ws.Range("E50").Value = "  +0.16%  "

$ws.Range("E51").Value = "  +3.51%  "
